$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 33 ("Ave flux um/m2"), shifting
# the rows below (Ave flux, Ave k, Ave k600, ER, GPP) down by one.
$ws.Rows.Item(33).Insert()

# Populate the new row with the average air pressure reading.
$ws.Range("A33").Value = "Ave Air Press (kPa)"
$ws.Range("B33").Value = 64.430000000000007

# The newly typed value cell keeps the default/"Normal" style (no special
# font), unlike the other numeric cells in column B of this block.
$ws.Range("B33").Style = "Normal"

# Match the author's final selection/cursor position.
[void]$ws.Range("B33").Select()
